# Rename the values with a dash rather then space
# (replace the space with an underscore in the "OS"/"PFS" mito/caba labels,
#  both the plain labels and the "Model: ..." labels)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("exponential")
$ws1.Range("B1").Value = "OS_mito"
$ws1.Range("C1").Value = "OS_caba"
$ws1.Range("D1").Value = "PFS_mito"
$ws1.Range("E1").Value = "PFS_caba"

$modelSheets = @("weibull", "lognormal", "loglogistic")
foreach ($sheetName in $modelSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A1").Value = "Model: PFS_caba"
    $ws.Range("A6").Value = "Model: PFS_mito"
    $ws.Range("A11").Value = "Model: OS_caba"
    $ws.Range("A16").Value = "Model: OS_mito"
}
